# Realestate Update resale numbers 2024-01-11 22:15
# Append one new data row (row 49) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 49

# Date / Week are text-like values ("2024-01-11", "01") that Excel would
# otherwise auto-convert to a date serial / number. A leading apostrophe
# forces them to stay plain text, same as a user typing them in manually.
$ws.Cells.Item($row, 1).Value  = "'2024-01-11"
$ws.Cells.Item($row, 2).Value  = "22:15:35"
$ws.Cells.Item($row, 3).Value  = "Thursday"
$ws.Cells.Item($row, 4).Value  = "'01"
$ws.Cells.Item($row, 5).Value  = 139555
$ws.Cells.Item($row, 6).Value  = 142801
$ws.Cells.Item($row, 7).Value  = 171928
$ws.Cells.Item($row, 8).Value  = 148234
$ws.Cells.Item($row, 9).Value  = -1
$ws.Cells.Item($row, 10).Value = 119677
$ws.Cells.Item($row, 11).Value = 224913
$ws.Cells.Item($row, 12).Value = 252438
$ws.Cells.Item($row, 13).Value = 185168
$ws.Cells.Item($row, 14).Value = 110498
$ws.Cells.Item($row, 15).Value = 40815
$ws.Cells.Item($row, 16).Value = 30906
$ws.Cells.Item($row, 17).Value = 72949
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42545
$ws.Cells.Item($row, 20).Value = -1
